# Applies the IFRS financial-data correction described in the commit
# "error solve ifrs list": rows 2-6 get their D:AJ figures replaced with
# corrected (much smaller / different) numbers, row 6 additionally loses
# its AG/AH entries, and rows 7-9 lose all data beyond A:C (D:AI cleared).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update financial figures
$ws.Range("D2").Value = 6723
$ws.Range("E2").Value = 19
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = -33
$ws.Range("H2").Value = -48
$ws.Range("I2").Value = -64
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 5164
$ws.Range("L2").Value = 3430
$ws.Range("M2").Value = 1734
$ws.Range("N2").Value = 1639
$ws.Range("O2").Value = 96
$ws.Range("P2").Value = 132
$ws.Range("Q2").Value = 151
$ws.Range("R2").Value = -487
$ws.Range("S2").Value = 267
$ws.Range("T2").Value = 462
$ws.Range("U2").Value = -310
$ws.Range("V2").Value = 1409
$ws.Range("W2").Value = 0.28
$ws.Range("X2").Value = -0.71
$ws.Range("Y2").Value = -3.79
$ws.Range("Z2").Value = -0.96
$ws.Range("AA2").Value = 197.73
$ws.Range("AB2").Value = 1205.36
$ws.Range("AC2").Value = -204
$ws.Range("AD2").Value = -9.76
$ws.Range("AE2").Value = 5250
$ws.Range("AF2").Value = 0.38
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 31212528

# Row 3: update financial figures
$ws.Range("D3").Value = 6223
$ws.Range("E3").Value = 128
$ws.Range("F3").Value = 128
$ws.Range("G3").Value = 95
$ws.Range("H3").Value = 49
$ws.Range("I3").Value = 46
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 5707
$ws.Range("L3").Value = 3883
$ws.Range("M3").Value = 1824
$ws.Range("N3").Value = 1702
$ws.Range("O3").Value = 122
$ws.Range("P3").Value = 132
$ws.Range("Q3").Value = 377
$ws.Range("R3").Value = -699
$ws.Range("S3").Value = 498
$ws.Range("T3").Value = 706
$ws.Range("U3").Value = -329
$ws.Range("V3").Value = 1903
$ws.Range("W3").Value = 2.07
$ws.Range("X3").Value = 0.78
$ws.Range("Y3").Value = 2.74
$ws.Range("Z3").Value = 0.9
$ws.Range("AA3").Value = 212.89
$ws.Range("AB3").Value = 1247.34
$ws.Range("AC3").Value = 147
$ws.Range("AD3").Value = 14.8
$ws.Range("AE3").Value = 5452
$ws.Range("AF3").Value = 0.4
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 31212528

# Row 4: update financial figures
$ws.Range("D4").Value = 6260
$ws.Range("E4").Value = 182
$ws.Range("F4").Value = 182
$ws.Range("G4").Value = 207
$ws.Range("H4").Value = 156
$ws.Range("I4").Value = 151
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 6310
$ws.Range("L4").Value = 4393
$ws.Range("M4").Value = 1917
$ws.Range("N4").Value = 1790
$ws.Range("O4").Value = 127
$ws.Range("P4").Value = 132
$ws.Range("Q4").Value = 62
$ws.Range("R4").Value = -652
$ws.Range("S4").Value = 460
$ws.Range("T4").Value = 666
$ws.Range("U4").Value = -604
$ws.Range("V4").Value = 2338
$ws.Range("W4").Value = 2.9
$ws.Range("X4").Value = 2.5
$ws.Range("Y4").Value = 8.63
$ws.Range("Z4").Value = 2.6
$ws.Range("AA4").Value = 229.18
$ws.Range("AB4").Value = 1358.16
$ws.Range("AC4").Value = 483
$ws.Range("AD4").Value = 4.96
$ws.Range("AE4").Value = 5734
$ws.Range("AF4").Value = 0.42
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 31212528

# Row 5: update financial figures
$ws.Range("D5").Value = 6020
$ws.Range("E5").Value = -367
$ws.Range("F5").Value = -367
$ws.Range("G5").Value = -501
$ws.Range("H5").Value = -411
$ws.Range("I5").Value = -338
$ws.Range("J5").Value = -73
$ws.Range("K5").Value = 5491
$ws.Range("L5").Value = 3674
$ws.Range("M5").Value = 1817
$ws.Range("N5").Value = 1758
$ws.Range("O5").Value = 59
$ws.Range("P5").Value = 267
$ws.Range("Q5").Value = 713
$ws.Range("R5").Value = -504
$ws.Range("S5").Value = -274
$ws.Range("T5").Value = 618
$ws.Range("U5").Value = 96
$ws.Range("V5").Value = 1741
$ws.Range("W5").Value = -6.1
$ws.Range("X5").Value = -6.83
$ws.Range("Y5").Value = -19.04
$ws.Range("Z5").Value = -6.97
$ws.Range("AA5").Value = 202.17
$ws.Range("AB5").Value = 609.25
$ws.Range("AC5").Value = -951
$ws.Range("AD5").Value = -1.12
$ws.Range("AE5").Value = 3298
$ws.Range("AF5").Value = 0.32
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 53312072

# Row 6: update financial figures
$ws.Range("D6").Value = 6509
$ws.Range("E6").Value = -121
$ws.Range("F6").Value = -121
$ws.Range("G6").Value = -110
$ws.Range("H6").Value = -101
$ws.Range("I6").Value = -96
$ws.Range("K6").Value = 5703
$ws.Range("L6").Value = 3909
$ws.Range("M6").Value = 1793
$ws.Range("N6").Value = 1741
$ws.Range("P6").Value = 267
$ws.Range("Q6").Value = 507
$ws.Range("R6").Value = -228
$ws.Range("S6").Value = -36
$ws.Range("T6").Value = 504
$ws.Range("U6").Value = 3
$ws.Range("V6").Value = 1708
$ws.Range("W6").Value = -1.87
$ws.Range("X6").Value = -1.55
$ws.Range("Y6").Value = -5.48
$ws.Range("Z6").Value = -1.81
$ws.Range("AA6").Value = 217.98
$ws.Range("AB6").Value = 611.47
$ws.Range("AC6").Value = -180
$ws.Range("AD6").Value = -4.65
$ws.Range("AE6").Value = 3265
$ws.Range("AF6").Value = 0.26
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 53312072

# Row 6: AG/AH no longer reported - remove them entirely
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9: all figures beyond A:C removed (no longer reported)
$ws.Range("D7:AI9").ClearContents()

